$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08815322733340154
$ws.Range("H2").Value = 32.92536994944349
$ws.Range("I2").Value = -8.467553196120466
$ws.Range("G3").Value = 0.09720123653599541
$ws.Range("H3").Value = -17.80917749778218
$ws.Range("G4").Value = -0.39586394350955
$ws.Range("H4").Value = -43.95295855841973
$ws.Range("G5").Value = -0.3385672669325922
$ws.Range("H5").Value = 15.14920931535176
$ws.Range("G6").Value = 0.224871422970643
$ws.Range("H6").Value = 14.06155540193005
$ws.Range("G7").Value = 0.272224770871786
$ws.Range("H7").Value = 31.26713205022098
$ws.Range("G8").Value = 0.08645905805581128
$ws.Range("H8").Value = -15.14971802964202
$ws.Range("G9").Value = 0.1269090409631717
$ws.Range("H9").Value = 0.3391678660798939
$ws.Range("G10").Value = -0.0002268147440498427
$ws.Range("H10").Value = -100.3691778133202
$ws.Range("G11").Value = 0.03396086582882882
$ws.Range("H11").Value = -31.98325043328909
$ws.Range("G12").Value = 0.1371822776914079
$ws.Range("H12").Value = 48.19748355984777
$ws.Range("G13").Value = 0.09569972464618243
$ws.Range("H13").Value = 25.57746948261492
$ws.Range("G14").Value = 0.2167102398658949
$ws.Range("H14").Value = -4.094236132987986
$ws.Range("G15").Value = 0.2640720401677575
$ws.Range("H15").Value = 7.188379130012581
$ws.Range("G16").Value = 0.1429041827411249
$ws.Range("H16").Value = 25.63573009145004
$ws.Range("G17").Value = 0.0943169704345611
$ws.Range("H17").Value = -36.87304064731713
$ws.Range("G18").Value = -0.002855643521567395
$ws.Range("H18").Value = 68.10031569503087
$ws.Range("G19").Value = 0.001865998991759696
$ws.Range("H19").Value = -92.2961967212396
$ws.Range("G20").Value = 0.09873453375979419
$ws.Range("H20").Value = 16.07475656511287
$ws.Range("G21").Value = 0.1450757409329498
$ws.Range("H21").Value = 121.6454545962644
$ws.Range("G22").Value = 0.189675963253995
$ws.Range("H22").Value = -0.9843907693012284
$ws.Range("G23").Value = 0.1963916179169492
$ws.Range("H23").Value = -8.954589352080738
$ws.Range("G24").Value = -0.009526415772110768
$ws.Range("H24").Value = -150.3283132985683
$ws.Range("G25").Value = -0.03330076571945329
$ws.Range("H25").Value = -43.18401244063176
$ws.Range("G26").Value = 0.1917716490241828
$ws.Range("H26").Value = -6.391813082429517
$ws.Range("G27").Value = 0.2238005697415886
$ws.Range("H27").Value = 16.02825670751736
$ws.Range("G28").Value = 0.02228585007918493
$ws.Range("H28").Value = -66.69462451552401
$ws.Range("G29").Value = 0.05635487996235866
$ws.Range("H29").Value = -40.2162706097165
